$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.960.08'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.890.09'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.016'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.015'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4701'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3923'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.96'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07979'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.015'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '1.889.88'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.970'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.151'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06776'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001049'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = '27.947.38'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.497'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.367'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.82%  '
$ws.Range('D26').Value = '2.108.46'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.101'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.506'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09576'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9629'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.651'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.352'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.363'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06138'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02250'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.213'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.196'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5957'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1900'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.270'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5689'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.947'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.402'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06858'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '114.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.070'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.66%  '
